$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A:E -> B:F)
$ws.Range("A1:A8").EntireColumn.Insert()

# New column A: Roll No header + sequential numbers
$ws.Range("A1").Value = "Roll No"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7

# Copy style of header row (B1, now bold/centered/bordered) onto new A1
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122) # xlPasteFormats

# Replace the old last data column (now column F, previously E = 2/14/23) with new date/header and values
$ws.Range("F1").Value = "'4/5/23"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats, keep same header style as other date columns
$ws.Range("F2").Value = "PRESENT"
$ws.Range("F3").Value = "ABSENT"
$ws.Range("F4").Value = "ABSENT"
$ws.Range("F5").Value = "PRESENT"
$ws.Range("F6").Value = "ABSENT"
$ws.Range("F7").Value = "ABSENT"
$ws.Range("F8").Value = "ABSENT"

$ws.Range("A1").Select()
